$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 168, pushing the existing
# rows 168-178 down to 170-180 (same as Excel's Insert > Entire Row).
$ws.Range("A168:A169").EntireRow.Insert()

# New row 168: weekly price record for "Región de La Araucanía"
$ws.Range("A168").Value = 10
$ws.Range("B168").Value = "Vega Modelo de Temuco"
$ws.Range("C168").Value = "La Araucanía"
$ws.Range("D168").Value = 44578
$ws.Range("E168").Value = 9
$ws.Range("F168").Value = 100112052
$ws.Range("G168").Value = "Albahaca"
$ws.Range("H168").Value = "Sin especificar"
$ws.Range("I168").Value = "Primera"
$ws.Range("J168").Value = 30
$ws.Range("K168").Value = 5000
$ws.Range("L168").Value = 5000
$ws.Range("M168").Value = 5000
$ws.Range("N168").Value = "$/paquete"
$ws.Range("O168").Value = "Región de La Araucanía"
$ws.Range("P168").Value = 5000
$ws.Range("Q168").Value = 1
$ws.Range("R168").Value = "Hortaliza"

# New row 169: weekly price record for "Región del Maule"
$ws.Range("A169").Value = 10
$ws.Range("B169").Value = "Vega Modelo de Temuco"
$ws.Range("C169").Value = "La Araucanía"
$ws.Range("D169").Value = 44578
$ws.Range("E169").Value = 9
$ws.Range("F169").Value = 100112052
$ws.Range("G169").Value = "Albahaca"
$ws.Range("H169").Value = "Sin especificar"
$ws.Range("I169").Value = "Primera"
$ws.Range("J169").Value = 50
$ws.Range("K169").Value = 5000
$ws.Range("L169").Value = 5000
$ws.Range("M169").Value = 5000
$ws.Range("N169").Value = "$/paquete"
$ws.Range("O169").Value = "Región del Maule"
$ws.Range("P169").Value = 5000
$ws.Range("Q169").Value = 1
$ws.Range("R169").Value = "Hortaliza"
